# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" stat (header in G1 is "K"). Update the computed
# values for rows 2-13 to reflect the regenerated K values (replacing the
# old Strike# derived numbers).
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    7  = 0
    9  = 1
    10 = 2
    11 = 2
    12 = 0
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
